# Add MIAPPE example values: collapse the three example rows (2-4) of the
# "nutrients" sheet into a single example-values row (row 2), and update a
# few of the unit / term values to use ontology short names + PURL ids.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("nutrients")

# Remove the two extra example rows (rows 3 and 4); this also shrinks the
# worksheet dimension and the annotationTable/autoFilter ranges from
# A1:AJ4 down to A1:AJ2.
$ws.Range("A3:AJ4").EntireRow.Delete()

# Update the remaining example-values row (row 2) with the new values.
$ws.Range("E2").Value = ""
$ws.Range("S2").Value = ""
$ws.Range("B2").Value = "Ca (XEO:00058): 5 mg/L"
$ws.Range("F2").Value = "milligram per square meter"
$ws.Range("G2").Value = "UO"
$ws.Range("H2").Value = "http://purl.obolibrary.org/obo/UO_0000309"
$ws.Range("I2").Value = "nitrogen: [concentration]; phosphorus: [concentration]"
$ws.Range("L2").Value = "Ca (XEO:00058): 5 mg/L"
$ws.Range("O2").Value = "[mg/m2]"
$ws.Range("T2").Value = "liter"
$ws.Range("U2").Value = "UO"
$ws.Range("V2").Value = "http://purl.obolibrary.org/obo/UO_0000099"
$ws.Range("W2").Value = "-10 to -30 kPa"
$ws.Range("Z2").Value = "drip irrigation"
$ws.Range("AC2").Value = "Ca (XEO:00058): 5 mg/L"
$ws.Range("AF2").Value = "dS m-1"
